$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bcde = New-Object 'object[,]' 21,4
$g = New-Object 'object[,]' 21,1

$bcde[0,0] = 3.182878228561681; $bcde[0,1] = 1.65323645889881; $bcde[0,2] = 0.7127328510149897; $bcde[0,3] = 0.4998867070740569
$g[0,0] = 6.048734245549538
$bcde[1,0] = 0.02258322285507441; $bcde[1,1] = 0.004309184025731883; $bcde[1,2] = 3.082599426703578; $bcde[1,3] = 0.4998867070740569
$g[1,0] = 3.609378540658442
$bcde[2,0] = 3.182878228561681; $bcde[2,1] = 1.65323645889881; $bcde[2,2] = 0.1529057820181812; $bcde[2,3] = 0.4998867070740569
$g[2,0] = 5.488907176552729
$bcde[3,0] = 3.182878228561681; $bcde[3,1] = 1.65323645889881; $bcde[3,2] = 0.7127328510149897; $bcde[3,3] = 0.4998867070740569
$g[3,0] = 6.048734245549538
$bcde[4,0] = 3.182878228561681; $bcde[4,1] = 1.65323645889881; $bcde[4,2] = 0.7127328510149897; $bcde[4,3] = 0.4998867070740569
$g[4,0] = 6.048734245549538
$bcde[5,0] = 3.182878228561681; $bcde[5,1] = 1.65323645889881; $bcde[5,2] = 0.7127328510149897; $bcde[5,3] = 0.4998867070740569
$g[5,0] = 6.048734245549538
$bcde[6,0] = 0.06328177979961902; $bcde[6,1] = 0.004309184025731883; $bcde[6,2] = 0.7127328510149897; $bcde[6,3] = 0.4998867070740569
$g[6,0] = 1.280210521914398
$bcde[7,0] = 0.3464964993005633; $bcde[7,1] = 1.65323645889881; $bcde[7,2] = 3.082599426703578; $bcde[7,3] = 6.48142807727062
$g[7,0] = 11.56376046217357
$bcde[8,0] = 3.182878228561681; $bcde[8,1] = 1.65323645889881; $bcde[8,2] = 0.7127328510149897; $bcde[8,3] = 0.4998867070740569
$g[8,0] = 6.048734245549538
$bcde[9,0] = 3.182878228561681; $bcde[9,1] = 9.226618575922256; $bcde[9,2] = 157.8057217802531; $bcde[9,3] = 6.48142807727062
$g[9,0] = 176.6966466620077
$bcde[10,0] = 1.505614041169197; $bcde[10,1] = 0.3375848360084654; $bcde[10,2] = 0.1529057820181812; $bcde[10,3] = 0.4998867070740569
$g[10,0] = 2.495991366269901
$bcde[11,0] = 3.182878228561681; $bcde[11,1] = 1.65323645889881; $bcde[11,2] = 0.7127328510149897; $bcde[11,3] = 0.4998867070740569
$g[11,0] = 6.048734245549538
$bcde[12,0] = 3.182878228561681; $bcde[12,1] = 1.65323645889881; $bcde[12,2] = 0.7127328510149897; $bcde[12,3] = 0.4998867070740569
$g[12,0] = 6.048734245549538
$bcde[13,0] = 1.505614041169197; $bcde[13,1] = 9.226618575922256; $bcde[13,2] = 0.1529057820181812; $bcde[13,3] = 6.48142807727062
$g[13,0] = 17.36656647638025
$bcde[14,0] = 0.1554434735375247; $bcde[14,1] = 0.3375848360084654; $bcde[14,2] = 0.7127328510149897; $bcde[14,3] = 0.4998867070740569
$g[14,0] = 1.705647867635037
$bcde[15,0] = 3.182878228561681; $bcde[15,1] = 1.65323645889881; $bcde[15,2] = 0.7127328510149897; $bcde[15,3] = 0.4998867070740569
$g[15,0] = 6.048734245549538
$bcde[16,0] = 0.001754667048134761; $bcde[16,1] = 0.3375848360084654; $bcde[16,2] = 0.1529057820181812; $bcde[16,3] = 0.4998867070740569
$g[16,0] = 0.9921319921488383
$bcde[17,0] = 3.182878228561681; $bcde[17,1] = 1.65323645889881; $bcde[17,2] = 0.7127328510149897; $bcde[17,3] = 0.4998867070740569
$g[17,0] = 6.048734245549538
$bcde[18,0] = 3.182878228561681; $bcde[18,1] = 1.65323645889881; $bcde[18,2] = 0.1529057820181812; $bcde[18,3] = 6.48142807727062
$g[18,0] = 11.47044854674929
$bcde[19,0] = 3.182878228561681; $bcde[19,1] = 1.65323645889881; $bcde[19,2] = 3.082599426703578; $bcde[19,3] = 0.4998867070740569
$g[19,0] = 8.418600821238126
$bcde[20,0] = 0.1554434735375247; $bcde[20,1] = 0.05231270169004087; $bcde[20,2] = 3.082599426703578; $bcde[20,3] = 0.4998867070740569
$g[20,0] = 3.790242309005201

$ws.Range("B2:E22").Value = $bcde
$ws.Range("G2:G22").Value = $g
